$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 11 (ATP8A1) - content shifts up
$ws.Rows("11:11").Delete()

# Delete row (now) 24, which was originally row 25 (TULP3) - content shifts up
$ws.Rows("24:24").Delete()

# Update the title cell A1 text (count of genes changed from 17 to 15)
$ws.Range("A1").Value = "15 DE Custom Ciliome expressed in scRNA-seq epithelial clusters"
